# Mise a jour des donnees (add latest scrape rows to France / Monde / percent)

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd h:mm:ss"

# ---------------------------------------------------------------
# Sheet "France": two new rows (15, 16) in columns A:D
# ---------------------------------------------------------------
$wsFrance = $wb.Worksheets.Item("France")

$wsFrance.Cells.Item(15, 1).Value = 43912.69236697095
$wsFrance.Cells.Item(15, 1).NumberFormat = $dateFmt
$wsFrance.Cells.Item(15, 2).Value = 14485
$wsFrance.Cells.Item(15, 3).Value = 562
$wsFrance.Cells.Item(15, 4).Value = 12

$wsFrance.Cells.Item(16, 1).Value = 43912.76475450717
$wsFrance.Cells.Item(16, 1).NumberFormat = $dateFmt
$wsFrance.Cells.Item(16, 2).Value = 14485
$wsFrance.Cells.Item(16, 3).Value = 562
$wsFrance.Cells.Item(16, 4).Value = 12

# ---------------------------------------------------------------
# Sheet "Monde": two new rows (10, 11) in columns A:D
# ---------------------------------------------------------------
$wsMonde = $wb.Worksheets.Item("Monde")

$wsMonde.Cells.Item(10, 1).Value = 43912.69236697095
$wsMonde.Cells.Item(10, 1).NumberFormat = $dateFmt
$wsMonde.Cells.Item(10, 2).Value = 316652
$wsMonde.Cells.Item(10, 3).Value = 13598
$wsMonde.Cells.Item(10, 4).Value = 94176

$wsMonde.Cells.Item(11, 1).Value = 43912.76475450717
$wsMonde.Cells.Item(11, 1).NumberFormat = $dateFmt
$wsMonde.Cells.Item(11, 2).Value = 318662
$wsMonde.Cells.Item(11, 3).Value = 13672
$wsMonde.Cells.Item(11, 4).Value = 94704

# ---------------------------------------------------------------
# Sheet "percent": twelve new rows (47-58) in columns A:E
# ---------------------------------------------------------------
$wsPercent = $wb.Worksheets.Item("percent")

$batch1Date = 43912.69242277653
$batch1 = @(
    @("France",          4.57, 4.13, 0.01),
    @("Italie",          16.92, 35.48, 6.45),
    @("Espagne",         9.02, 12.65, 2.26),
    @("Allemagne",       7.55, 0.68, 0.28),
    @("UK",               1.6, 1.72, 0.07000000000000001),
    @("Reste du monde",  60.34, 45.34, 90.93000000000001)
)

$row = 47
foreach ($entry in $batch1) {
    $wsPercent.Cells.Item($row, 1).Value = $batch1Date
    $wsPercent.Cells.Item($row, 1).NumberFormat = $dateFmt
    $wsPercent.Cells.Item($row, 2).Value = $entry[0]
    $wsPercent.Cells.Item($row, 3).Value = $entry[1]
    $wsPercent.Cells.Item($row, 4).Value = $entry[2]
    $wsPercent.Cells.Item($row, 5).Value = $entry[3]
    $row = $row + 1
}

$batch2Date = 43912.76479431608
$batch2 = @(
    @("France",          4.55, 4.11, 0.01),
    @("Italie",          16.81, 35.29, 6.41),
    @("Espagne",         8.98, 12.84, 2.72),
    @("Allemagne",       7.52, 0.67, 0.28),
    @("UK",               1.59, 1.71, 0.07000000000000001),
    @("Reste du monde",  60.55, 45.38, 90.51000000000001)
)

foreach ($entry in $batch2) {
    $wsPercent.Cells.Item($row, 1).Value = $batch2Date
    $wsPercent.Cells.Item($row, 1).NumberFormat = $dateFmt
    $wsPercent.Cells.Item($row, 2).Value = $entry[0]
    $wsPercent.Cells.Item($row, 3).Value = $entry[1]
    $wsPercent.Cells.Item($row, 4).Value = $entry[2]
    $wsPercent.Cells.Item($row, 5).Value = $entry[3]
    $row = $row + 1
}
